# Update "want to go" counts (column F) on the "展览" and "全部类型" sheets
# to reflect newly generated output data (gh-pages refresh).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 185
$ws1.Range("F4").Value  = 5257
$ws1.Range("F7").Value  = 21
$ws1.Range("F8").Value  = 591
$ws1.Range("F9").Value  = 551
$ws1.Range("F12").Value = 1451
$ws1.Range("F13").Value = 4198
$ws1.Range("F15").Value = 175
$ws1.Range("F16").Value = 158
$ws1.Range("F17").Value = 95
$ws1.Range("F18").Value = 3285
$ws1.Range("F23").Value = 190
$ws1.Range("F24").Value = 114
$ws1.Range("F26").Value = 136
$ws1.Range("F28").Value = 297
$ws1.Range("F31").Value = 15
$ws1.Range("F32").Value = 19
$ws1.Range("F33").Value = 20

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 185
$ws4.Range("F5").Value  = 5257
$ws4.Range("F8").Value  = 21
$ws4.Range("F9").Value  = 591
$ws4.Range("F10").Value = 551
$ws4.Range("F13").Value = 1451
$ws4.Range("F14").Value = 4198
$ws4.Range("F16").Value = 175
$ws4.Range("F17").Value = 158
$ws4.Range("F18").Value = 95
$ws4.Range("F19").Value = 3285
$ws4.Range("F24").Value = 190
$ws4.Range("F25").Value = 114
$ws4.Range("F27").Value = 136
$ws4.Range("F29").Value = 297
$ws4.Range("F32").Value = 15
$ws4.Range("F33").Value = 19
$ws4.Range("F34").Value = 20
